# Switch off automatic recalculation first, so that the cached results of
# existing formulas (e.g. the shared "=SUM(C+D)" formula in column B) are
# left untouched when the value of D7 stops being numeric.
$excel.Calculation = -4135  # xlCalculationManual

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell D7 currently holds the number 548; the source data now records it
# as the time-like text "5:48" (same style as the neighbouring D6 value
# "4:36"), so replace the numeric value with that text string.
$ws.Range("D7").Value = "5:48"
